# Generate DB View SQL files for Dataprocessor all tables
# Adds a TEMPLATE column (B) entry for the cre_table / cre_view rows and
# renames the dataprocessor "all" view script to include the .sql suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "TEMPLATE" column (B) values for the cre_table / cre_view rows ---
$ws.Range("B36").Value = "template_cre_table.sql"
$ws.Range("B38").Value = "template_cre_table.sql"
$ws.Range("B40").Value = "template_cre_table.sql"
$ws.Range("B42").Value = "template_cre_table.sql"
$ws.Range("B45").Value = "template_cre_view.sql"
$ws.Range("B47").Value = "template_cre_view2.sql"

# --- Row 47: dataprocessor "all" view script gets a .sql extension, and
#     picks up an additional OWNER_SCHEMA (L) value ---
$ws.Range("A47").Value = "19_cre_view_typ_dataproc_all.sql"
$ws.Range("L47").Value = "db_log"

# --- Column widths: widen the new TEMPLATE column (B) and split the old
#     shared 16.28515625-wide C:D range into two independently sized
#     columns now that their content differs in length. ---
$ws.Columns.Item(2).ColumnWidth = 33.666667
$ws.Columns.Item(3).ColumnWidth = 21
$ws.Columns.Item(4).ColumnWidth = 20

# --- Selection moves from the old last row (A48) to the edited row (A47) ---
$ws.Range("A47").Select()
